# Update the table values for section 1 (root-finding method iterations).
# Cells store their numbers as text (inline strings), so each write first
# forces Text format to stop Excel auto-converting the numeric-looking
# string to a real number, then resets the style back to Normal so the
# cell's style index is left untouched (matches the original workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 (A2=0 unchanged) - B2/C2 change, D2 unchanged
Set-TextValue "B2" "6.0"
Set-TextValue "C2" "34.0"

# Row 3
Set-TextValue "B3" "3.16666666666667"
Set-TextValue "C3" "8.02777777777778"
Set-TextValue "D3" "2.83333333333333"

# Row 4
Set-TextValue "B4" "1.89912280701754"
Set-TextValue "C4" "1.6066674361342"
Set-TextValue "D4" "1.26754385964912"

# Row 5
Set-TextValue "B5" "1.47612029496374"
Set-TextValue "C5" "0.178931125203831"
Set-TextValue "D5" "0.423002512053807"

# Row 6
Set-TextValue "B6" "1.41551170980496"
Set-TextValue "C6" "0.0036734005949488"
Set-TextValue "D6" "0.0606085851587816"

# Row 7
Set-TextValue "B7" "1.41421415763018"
Set-TextValue "C7" "1.68364164609969e-06"
Set-TextValue "D7" "0.0012975521747733"

# Row 8
Set-TextValue "B8" "1.41421356237322"
Set-TextValue "C8" "3.5438318946034997e-13"
Set-TextValue "D8" "5.95256961943846e-07"

# Remove rows 9, 10, 11 (the table now only spans 7 iterations, rows 2-8)
$ws.Range("A9:D11").EntireRow.Delete()
